$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp string (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 26 de Junio de 2020 a las 22:16"

# --- Country name swaps (rows whose shared-string slot swapped with its neighbor) ---
$ws.Range("A128").Value = "Yemen"
$ws.Range("A129").Value = "Congo"
$ws.Range("A200").Value = "Laos"
$ws.Range("A201").Value = "Santa Lucia"
$ws.Range("A202").Value = "Dominica"
$ws.Range("A203").Value = "Fiyi"
$ws.Range("A208").Value = "Islas Malvinas"
$ws.Range("A209").Value = "Groenlandia"
$ws.Range("A212").Value = "Montserrat"
$ws.Range("A213").Value = "Seychelles"

# --- Updated covid-19 numeric stats ---
$ws.Range("B4").Value = 2538867
$ws.Range("C4").Value = 34279
$ws.Range("D4").Value = 1056669
$ws.Range("E4").Value = 1354899
$ws.Range("G4").Value = 519
$ws.Range("H4").Value = 127299
$ws.Range("B7").Value = 509446
$ws.Range("C7").Value = 18276
$ws.Range("E7").Value = 197840
$ws.Range("B17").Value = 194256
$ws.Range("C17").Value = 471
$ws.Range("E17").Value = 8130
$ws.Range("G17").Value = 14
$ws.Range("H17").Value = 9026
$ws.Range("B27").Value = 62755
$ws.Range("C27").Value = 1625
$ws.Range("D27").Value = 16737
$ws.Range("E27").Value = 43398
$ws.Range("G27").Value = 87
$ws.Range("H27").Value = 2620
$ws.Range("B75").Value = 7427
$ws.Range("C75").Value = 250
$ws.Range("E75").Value = 2369
$ws.Range("B94").Value = 3907
$ws.Range("C94").Value = 168
$ws.Range("D94").Value = 1280
$ws.Range("E94").Value = 2507
$ws.Range("G94").Value = 1
$ws.Range("H94").Value = 120
$ws.Range("B100").Value = 2836
$ws.Range("C100").Value = 152
$ws.Range("D100").Value = 1280
$ws.Range("E100").Value = 1544
$ws.Range("B123").Value = 1394
$ws.Range("C123").Value = 40
$ws.Range("D123").Value = 914
$ws.Range("E123").Value = 421
$ws.Range("G123").Value = 3
$ws.Range("H123").Value = 59
$ws.Range("B128").Value = 1089
$ws.Range("C128").Value = 13
$ws.Range("D128").Value = 402
$ws.Range("E128").Value = 394
$ws.Range("G128").Value = 5
$ws.Range("H128").Value = 293
$ws.Range("B129").Value = 1087
$ws.Range("D129").Value = 456
$ws.Range("E129").Value = 594
$ws.Range("H129").Value = 37
$ws.Range("B152").Value = 561
$ws.Range("C152").Value = 10
$ws.Range("D152").Value = 135
$ws.Range("E152").Value = 420
$ws.Range("D163").Value = 102
$ws.Range("E163").Value = 144
$ws.Range("G163").Value = 1
$ws.Range("H163").Value = 8
$ws.Range("B178").Value = 121
$ws.Range("C178").Value = 19
$ws.Range("D178").Value = 22
$ws.Range("E178").Value = 99
$ws.Range("D212").Value = 10
$ws.Range("H212").Value = 1
$ws.Range("D213").Value = 11
$ws.Range("H213").Value = 0
